$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update odds values in row 5 per diff
$ws.Range("G5").Value = 2.3
$ws.Range("I5").Value = 3.3
$ws.Range("J5").Value = 3.1
$ws.Range("L5").Value = 4
$ws.Range("M5").Value = 1.11
$ws.Range("N5").Value = 6.5
$ws.Range("X5").Value = 10
$ws.Range("Z5").Value = 21
$ws.Range("AB5").Value = 41
$ws.Range("AC5").Value = 6.5
$ws.Range("AF5").Value = 67
$ws.Range("AL5").Value = 34
$ws.Range("AO5").Value = 13
